$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Cason Wallace", "PG,SG", "Oklahoma City Thunder"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Jalen Johnson", "SF,PF", "Atlanta Hawks"),
    @("Mark Williams", "C", "Charlotte Hornets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
